$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.880.91'
$ws.Range('E2').Value = '  +0.52%  '

# Row 3
$ws.Range('D3').Value = '1.895.13'
$ws.Range('E3').Value = '  +0.52%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').Value = '''0.7844'
$ws.Range('E5').Value = '  -0.92%  '

# Row 6
$ws.Range('D6').Value = '''243.88'
$ws.Range('E6').Value = '  +1.13%  '

# Row 7
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('D8').Value = '''0.3141'
$ws.Range('E8').Value = '  -0.69%  '

# Row 9
$ws.Range('D9').Value = '''25.77'
$ws.Range('E9').Value = '  +1.73%  '

# Row 10
$ws.Range('D10').Value = '''0.07272'
$ws.Range('E10').Value = '  +4.15%  '

# Row 11
$ws.Range('E11').Value = '  +0.89%  '

# Row 12
$ws.Range('D12').Value = '''0.7764'
$ws.Range('E12').Value = '  +2.05%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.954.13'
$ws.Range('E13').Value = '  +3.54%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''5.484'
$ws.Range('E14').Value = '  +3.66%  '

# Row 15
$ws.Range('D15').Value = '''94.37'
$ws.Range('E15').Value = '  +2.48%  '

# Row 16
$ws.Range('D16').Value = '''6.209'
$ws.Range('E16').Value = '  +4.70%  '

# Row 17
$ws.Range('D17').Value = '29.832.17'
$ws.Range('E17').Value = '  +0.43%  '

# Row 18
$ws.Range('D18').Value = '''13.96'
$ws.Range('E18').Value = '  +1.01%  '

# Row 19
$ws.Range('D19').Value = '''246.33'
$ws.Range('E19').Value = '  +1.38%  '

# Row 20
$ws.Range('D20').Value = '''0.000007832'
$ws.Range('E20').Value = '  +2.15%  '

# Row 21
$ws.Range('D21').Value = '''1.000'
$ws.Range('E21').Value = '  +0.04%  '

# Row 22
$ws.Range('D22').Value = '''8.121'
$ws.Range('E22').Value = '  -0.51%  '

# Row 23
$ws.Range('D23').Value = '2.117.07'
$ws.Range('E23').Value = '  +1.42%  '

# Row 24
$ws.Range('D24').Value = '''1.001'
$ws.Range('E24').Value = '  -0.08%  '

# Row 25
$ws.Range('D25').Value = '''0.1592'
$ws.Range('E25').Value = '  -4.76%  '

# Row 26
$ws.Range('D26').Value = '''9.464'
$ws.Range('E26').Value = '  +2.05%  '

# Row 27
$ws.Range('D27').Value = '''164.11'
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('E28').Value = '  +1.04%  '

# Row 29
$ws.Range('D29').Value = '''2.025'
$ws.Range('E29').Value = '  -0.94%  '

# Row 30
$ws.Range('D30').Value = '''1.438'
$ws.Range('E30').Value = '  +3.52%  '

# Row 31
$ws.Range('D31').Value = '''1.544'
$ws.Range('E31').Value = '  +0.82%  '

# Row 32
$ws.Range('D32').Value = '''4.476'
$ws.Range('E32').Value = '  +2.48%  '

# Row 33
$ws.Range('D33').Value = '''0.05574'
$ws.Range('E33').Value = '  -1.58%  '

# Row 34
$ws.Range('D34').Value = '''4.078'
$ws.Range('E34').Value = '  +0.81%  '

# Row 35
$ws.Range('E35').Value = '  -1.09%  '

# Row 36
$ws.Range('D36').Value = '''0.7538'
$ws.Range('E36').Value = '  +2.87%  '

# Row 37
$ws.Range('E37').Value = '  +0.67%  '

# Row 38
$ws.Range('D38').Value = '''2.678'
$ws.Range('E38').Value = '  +1.96%  '

# Row 39
$ws.Range('E39').Value = '  +1.48%  '

# Row 40
$ws.Range('D40').Value = '''2.795'
$ws.Range('E40').Value = '  +1.02%  '

# Row 41
$ws.Range('D41').Value = '1.144.39'
$ws.Range('E41').Value = '  +12.21%  '

# Row 42
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '''74.30'
$ws.Range('E42').Value = '  +2.89%  '

# Row 43
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '''0.4463'
$ws.Range('E43').Value = '  +1.60%  '

# Row 44
$ws.Range('D44').Value = '''5.958'
$ws.Range('E44').Value = '  +2.55%  '

# Row 45
$ws.Range('D45').Value = '''0.8538'
$ws.Range('E45').Value = '  +2.24%  '

# Row 47
$ws.Range('E47').Value = '  +1.93%  '

# Row 48
$ws.Range('D48').Value = '''3.154'
$ws.Range('E48').Value = '  +8.64%  '

# Row 49
$ws.Range('D49').Value = '''102.01'
$ws.Range('E49').Value = '  -0.41%  '

# Row 50
$ws.Range('D50').Value = '''7.549'
$ws.Range('E50').Value = '  +1.96%  '

# Row 51
$ws.Range('D51').Value = '''9.755'
$ws.Range('E51').Value = '  -1.70%  '
